$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove rows 9-18 (old entries being dropped)
$ws.Rows("9:18").Delete()

# Remove all existing hyperlinks (Range.Hyperlinks.Delete() clears all on this runtime)
$ws.Range("F2").Hyperlinks.Delete()

# Row 2
$ws.Range("A2").Value = "2025-11-27 06:28:33"
$ws.Range("B2").Value = "GoogleAppSheetで行政書士向け案件管理アプリ開発"
$ws.Range("D2").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("G2").Value = 123
$ws.Range("H2").Value = "◆開発 ◇アプリ"
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5441932")
$ws.Range("F2").Style = "Hyperlink"

# Row 3
$ws.Range("A3").Value = "2025-11-27 06:28:33"
$ws.Range("B3").Value = "Javaプログラミング研修の演習サポート講師業務【経験不問】(再掲)"
$ws.Range("D3").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("G3").Value = 85
$ws.Range("H3").Value = "★Java"
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5442063")
$ws.Range("F3").Style = "Hyperlink"

# Row 4
$ws.Range("A4").Value = "2025-11-27 06:28:33"
$ws.Range("B4").Value = "クラウド(AWS/Azure) 運用管理 研修の演習サポート講師業務【経験不問】(再掲)"
$ws.Range("D4").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("G4").Value = 38
$ws.Range("H4").Value = "◇管理"
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5442064")
$ws.Range("F4").Style = "Hyperlink"

# Row 5
$ws.Range("A5").Value = "2025-11-27 06:28:33"
$ws.Range("B5").Value = "Access業務システムのクラウド化(ZOHO Creator使用)をお手伝いください!(再依頼)"
$ws.Range("D5").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("G5").Value = 40
$ws.Range("H5").ClearContents()
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5442153")
$ws.Range("F5").Style = "Hyperlink"

# Row 6
$ws.Range("A6").Value = "2025-11-27 06:28:33"
$ws.Range("B6").Value = "急募 限定公開 PR 限定公開の仕事"
$ws.Range("D6").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("G6").Value = 25
$ws.Range("H6").ClearContents()
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5440230")
$ws.Range("F6").Style = "Hyperlink"

# Row 7
$ws.Range("A7").Value = "2025-11-27 06:28:33"
$ws.Range("B7").Value = "【募集】Amazonフラットファイル(ブラウズノード検証)"
$ws.Range("D7").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("G7").Value = 13
$ws.Range("H7").ClearContents()
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5442106")
$ws.Range("F7").Style = "Hyperlink"

# Row 8
$ws.Range("A8").Value = "2025-11-27 06:28:33"
$ws.Range("B8").Value = "【急募】ex4ファイルをmq4ファイルに変換していただける方"
$ws.Range("D8").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("G8").Value = 10
$ws.Range("H8").ClearContents()
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5442169")
$ws.Range("F8").Style = "Hyperlink"

# Adjust column widths
$ws.Columns.Item(2).ColumnWidth = 51
$ws.Columns.Item(4).ColumnWidth = 28
$ws.Columns.Item(8).ColumnWidth = 12

